$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 146 (shifts existing rows 146-152 down to 147-153)
$ws.Rows.Item(146).Insert()

# Populate the new row 146 with data for a new weekly price observation,
# matching the surrounding rows' constant columns (A, B, C, E, F, G, H, I, N, O, Q, R)
$ws.Cells.Item(146, 1).Value = 8
$ws.Cells.Item(146, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(146, 3).Value = "Coquimbo"
$ws.Cells.Item(146, 4).Value = 44747
$ws.Cells.Item(146, 5).Value = 4
$ws.Cells.Item(146, 6).Value = 100112044
$ws.Cells.Item(146, 7).Value = "Perejil"
$ws.Cells.Item(146, 8).Value = "Sin especificar"
$ws.Cells.Item(146, 9).Value = "Primera"
$ws.Cells.Item(146, 10).Value = 2800
$ws.Cells.Item(146, 11).Value = 1500
$ws.Cells.Item(146, 12).Value = 2000
$ws.Cells.Item(146, 13).Value = 1750
$ws.Cells.Item(146, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(146, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(146, 16).Value = 1167
$ws.Cells.Item(146, 17).Value = 1.5
$ws.Cells.Item(146, 18).Value = "Hortaliza"
